# Insert a new data row at row 656 (pushing the existing rows 656:697 down to
# 657:698) and populate it with the new record: 2026/01/15, 木, 17, 201
#
# Row 655 already contains the literal text "2026/01/15" / "木" in columns A/B,
# so we copy that row down into the newly inserted row to preserve the
# original "text" cell typing/formatting (typing the date string directly
# into a General-formatted cell would get auto-converted to a date value),
# then overwrite just the two numeric columns with the new C/D values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A656").EntireRow.Insert()

$ws.Range("A655:D655").Copy()
$ws.Range("A656:D656").PasteSpecial()

$ws.Range("C656").Value = 17
$ws.Range("D656").Value = 201
